# Automatische test-sync: 2025-06-22 18:56:50
#
# Adds the new incoming mail-log entry ("Ruilen van product") as row 26 on
# the "Logs" sheet, extends the conditional-formatting ranges to cover it,
# and refreshes the "Dashboard" category-count table (rows 7-12) so it
# stays sorted by descending count now that "Retour / Terugbetaling" has
# moved from 1 to 2 occurrences.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append the new mail entry as row 26
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A26").Value = "Ruilen van product"
$logs.Range("B26").Value = "mailmind.test@zohomail.eu"
$logs.Range("C26").Value = "Kan ik dit product ruilen voor een andere maat?"
$logs.Range("D26").Value = "Retour / Terugbetaling"
$logs.Range("F26").Value = "2025-06-22 18:56:13"
$logs.Range("G26").Value = "Nee"

# Extend the two conditional-formatting blocks (Categorie / Beantwoord)
# from row 25 to row 26 so the new row is covered too.
$catRules = $logs.Range("D2:D25").FormatConditions
$catRules.Item(1).ModifyAppliesToRange($logs.Range("D2:D26"))

$answeredRules = $logs.Range("G2:G25").FormatConditions
$answeredRules.Item(1).ModifyAppliesToRange($logs.Range("G2:G26"))

# ---------------------------------------------------------------------
# 2. Dashboard sheet: re-sync the category/count table (rows 7-12)
#    "Retour / Terugbetaling" now has 2 hits (was 1), so it moves up,
#    pushing the others down one slot, matching a fresh sort-by-count.
# ---------------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A7").Value = "Retour / Terugbetaling"
$dash.Range("B7").Value = 2

$dash.Range("A8").Value = "Overig"
$dash.Range("B8").Value = 2

$dash.Range("A9").Value = "Offerte / Prijsaanvraag"
$dash.Range("B9").Value = 2

$dash.Range("A10").Value = "Klacht / Probleem"
$dash.Range("B10").Value = 1

$dash.Range("A12").Value = "Uitnodiging / Evenement"
$dash.Range("B12").Value = 1
